$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-02 Wednesday", "2025-07-03 Thursday"),
    @("514÷6=", "808÷9="),
    @("474÷3=", "777÷4="),
    @("666÷2=", "948÷4="),
    @("386÷3=", "110÷8="),
    @("671÷6=", "108÷9="),
    @("769÷4=", "876÷5="),
    @("749÷7=", "113÷7="),
    @("738÷2=", "314÷9="),
    @("459÷9=", "981÷9="),
    @("522÷3=", "551÷6="),
    @("821÷9=", "921÷9="),
    @("415÷3=", "273÷7="),
    @("483÷4=", "218÷4="),
    @("269÷7=", "552÷7="),
    @("298÷7=", "307÷4="),
    @("540÷6=", "460÷8="),
    @("523÷9=", "615÷2="),
    @("629÷7=", "931÷7="),
    @("782÷2=", "319÷6="),
    @("235÷8=", "191÷5="),
    @("757÷9=", "378÷8="),
    @("715÷5=", "124÷6="),
    @("245÷7=", "816÷5="),
    @("373÷4=", "426÷7="),
    @("857÷3=", "622÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
